# Revert "cant save conditional": restore the original Sheet1!B2:K11
# probability matrix (re-normalized row/column proportions) that the
# reverted commit had overwritten.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 10,10
$data[0,0] = [double]"6.8177905308464848E-2"
$data[0,1] = [double]"7.3542402826855119E-2"
$data[0,2] = [double]"0.15744843635436259"
$data[0,3] = [double]"7.534095423048065E-2"
$data[0,4] = [double]"0.24359414041448241"
$data[0,5] = [double]"1.235383579745194E-2"
$data[0,6] = [double]"9.823544903871477E-2"
$data[0,7] = [double]"2.6157269224507521E-2"
$data[0,8] = [double]"6.8170519512801397E-2"
$data[0,9] = [double]"5.1782942462224008E-2"
$data[1,0] = [double]"9.9139167862266855E-2"
$data[1,1] = [double]"3.222909305064782E-2"
$data[1,2] = [double]"0.128470783257356"
$data[1,3] = [double]"0.11652028070794899"
$data[1,4] = [double]"1.530220632665393E-2"
$data[1,5] = [double]"0.35084395023560799"
$data[1,6] = [double]"0.14063734527258359"
$data[1,7] = [double]"4.0429804091821793E-2"
$data[1,8] = [double]"5.3815560526969917E-2"
$data[1,9] = [double]"0.17271892125909419"
$data[2,0] = [double]"1.724533715925395E-2"
$data[2,1] = [double]"2.1157243816254421E-2"
$data[2,2] = [double]"0.11951289489623509"
$data[2,3] = [double]"0.12141942887407869"
$data[2,4] = [double]"2.454558816324531E-2"
$data[2,5] = [double]"4.5401281507890988E-2"
$data[2,6] = [double]"1.5801948907031869E-2"
$data[2,7] = [double]"8.6286427524827694E-3"
$data[2,8] = [double]"5.6466981522909937E-2"
$data[2,9] = [double]"0.14543361878715561"
$data[3,0] = [double]"3.7618364418938309E-2"
$data[3,1] = [double]"2.1054181389870441E-3"
$data[3,2] = [double]"1.721444738436036E-2"
$data[3,3] = [double]"1.191684689058569E-3"
$data[3,4] = [double]"1.0778501788251431E-3"
$data[3,5] = [double]"7.1181031688648436E-2"
$data[3,6] = [double]"4.8722675796681589E-2"
$data[3,7] = [double]"6.674987789656482E-3"
$data[3,8] = [double]"6.2142679592344019E-3"
$data[3,9] = [double]"1.2945735615556E-2"
$data[4,0] = [double]"1.9885222381635578E-2"
$data[4,1] = [double]"0.4123822143698469"
$data[4,2] = [double]"6.9176575600114767E-2"
$data[4,3] = [double]"0.13258595577525711"
$data[4,4] = [double]"2.6293011938007281E-2"
$data[4,5] = [double]"4.0714054202298733E-2"
$data[4,6] = [double]"2.2122728469844619E-2"
$data[4,7] = [double]"0.1751234601400119"
$data[4,8] = [double]"2.3489932885906041E-2"
$data[4,9] = [double]"0.14446216064960979"
$data[5,0] = [double]"0.65457675753228117"
$data[5,1] = [double]"0.1738221436984688"
$data[5,2] = [double]"3.659663999489942E-2"
$data[5,3] = [double]"0.16493798825969899"
$data[5,4] = [double]"0.56748811915143793"
$data[5,5] = [double]"0.13187813209005461"
$data[5,6] = [double]"0.51461680273900445"
$data[5,7] = [double]"0.11683942041569439"
$data[5,8] = [double]"0.19013588532604189"
$data[5,9] = [double]"6.4464694888229521E-2"
$data[6,0] = [double]"5.2223816355810616E-3"
$data[6,1] = [double]"1.5326855123674909E-2"
$data[6,2] = [double]"8.1290445981701681E-3"
$data[6,3] = [double]"1.354989627929558E-2"
$data[6,4] = [double]"1.5726813972857769E-2"
$data[6,5] = [double]"2.293749532523873E-3"
$data[6,6] = [double]"1.2641559125625489E-2"
$data[6,7] = [double]"4.8298692136538772E-3"
$data[6,8] = [double]"2.3821360510398542E-3"
$data[6,9] = [double]"6.4095118422858831E-3"
$data[7,0] = [double]"4.6714490674318508E-2"
$data[7,1] = [double]"8.2597173144876319E-3"
$data[7,2] = [double]"3.6787911632503417E-2"
$data[7,3] = [double]"1.5491900957761401E-2"
$data[7,4] = [double]"1.6820995214998451E-3"
$data[7,5] = [double]"7.3948490363758751E-2"
$data[7,6] = [double]"2.6336581511719779E-2"
$data[7,7] = [double]"7.6518152710696261E-3"
$data[7,8] = [double]"2.324136216753666E-2"
$data[7,9] = [double]"8.4041688225293815E-2"
$data[8,0] = [double]"1.893830703012913E-2"
$data[8,1] = [double]"0.2148409893992933"
$data[8,2] = [double]"7.1184927794956804E-2"
$data[8,3] = [double]"0.13452796045372289"
$data[8,4] = [double]"2.4643574543138501E-2"
$data[8,5] = [double]"4.593732080081777E-2"
$data[8,6] = [double]"1.7645509612852248E-2"
$data[8,7] = [double]"0.1056059043794432"
$data[8,8] = [double]"1.980280056342696E-2"
$data[8,9] = [double]"0.134314646843289"
$data[9,0] = [double]"3.1994261119081782E-2"
$data[9,1] = [double]"4.6186690223792702E-2"
$data[9,2] = [double]"0.35538270266823929"
$data[9,3] = [double]"0.22416913095290639"
$data[9,4] = [double]"7.9613933663220809E-2"
$data[9,5] = [double]"0.2249245805180882"
$data[9,6] = [double]"0.10323939952594149"
$data[9,7] = [double]"0.50800455852824655"
$data[9,8] = [double]"0.55628055348413286"
$data[9,9] = [double]"0.18332048615144189"

$ws.Range("B2:K11").Value = $data

# Matches the saved selection state (B12) recorded in the reverted file.
$ws.Range("B12").Select()
